# Edit Southampton_stats.xlsx per commit diff:
#  1. Rename the per-category stat sheets to their "spaced out" display names.
#  2. Bump every player's Age (YY-DDD) in column E (rows 4-37) forward by one day,
#     on every stats sheet.
#  3. Fix the "Unnamed: 4_level_0" / "Playing Time" header split on the two sheets
#     where it was mis-merged (Standard Stats + Playing Time): the "Playing Time"
#     label should live in F1 (merged F1:I1), not G1 (merged G1:I1).

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets -------------------------------------------------------
$wb.Worksheets.Item("StandardStats").Name      = "Standard Stats"
$wb.Worksheets.Item("ShootingStats").Name      = "Shooting Stats"
$wb.Worksheets.Item("PassingStats").Name       = "Passing Stats"
$wb.Worksheets.Item("PassTypes").Name          = "Pass Types"
$wb.Worksheets.Item("GoalShotCreation").Name   = "Goal & Shot Creation"
$wb.Worksheets.Item("DefensiveActions").Name   = "Defensive Actions"
$wb.Worksheets.Item("PlayingTime").Name        = "Playing Time"
$wb.Worksheets.Item("MiscStats").Name          = "Miscellaneous Stats"
# "Possession" is unchanged.

# --- 2. Bump the Age-in-days column (E4:E37) by one day on every stats sheet ---
$ages = @(
    "27-342",
    "20-254",
    "23-050",
    "28-343",
    "26-311",
    "28-243",
    "26-060",
    "19-009",
    "24-266",
    "28-039",
    "28-280",
    "23-243",
    "27-161",
    "20-360",
    "31-053",
    "30-297",
    "23-034",
    "22-294",
    "25-028",
    "35-108",
    "36-315",
    "25-337",
    "31-025",
    "30-034",
    "24-030",
    "23-100",
    "31-184",
    "23-302",
    "28-253",
    "28-175",
    "22-052",
    "18-246",
    "21-358",
    "22-111"
)

$statSheetNames = @(
    "Standard Stats",
    "Shooting Stats",
    "Passing Stats",
    "Pass Types",
    "Goal & Shot Creation",
    "Defensive Actions",
    "Possession",
    "Playing Time",
    "Miscellaneous Stats"
)

foreach ($sheetName in $statSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 4
    foreach ($age in $ages) {
        $ws.Cells.Item($row, 5).Value = $age
        $row = $row + 1
    }
}

# --- 3. Fix the Playing-Time header merge on "Standard Stats" and "Playing Time" ---
$fixSheetNames = @("Standard Stats", "Playing Time")

foreach ($sheetName in $fixSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Unmerge every row-1 merged block so we can rebuild them in the same
    # relative order as the target workbook (Y1:AH1, F1:I1, R1:U1, J1:Q1, V1:X1).
    $ws.Range("Y1:AH1").UnMerge()
    $ws.Range("R1:U1").UnMerge()
    $ws.Range("J1:Q1").UnMerge()
    $ws.Range("G1:I1").UnMerge()
    $ws.Range("V1:X1").UnMerge()

    # Move the "Playing Time" label from G1 to F1 (replacing the stray
    # "Unnamed: 4_level_0" placeholder), and blank out the old G1.
    $ws.Range("F1").Value = "Playing Time"
    $ws.Range("G1").Value = ""

    $ws.Range("Y1:AH1").Merge()
    $ws.Range("F1:I1").Merge()
    $ws.Range("R1:U1").Merge()
    $ws.Range("J1:Q1").Merge()
    $ws.Range("V1:X1").Merge()
}
